# Update latest output (run 130)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E4").Value = 570.2312459999999
$wsSchedule.Range("F4").Value = 30.17096539682539
$wsSchedule.Range("E5").Value = -189.16488825
$wsSchedule.Range("F5").Value = -5.560402358906526

# --- Detailed sheet updates (Price column B, Type column C) ---
$wsDetailed.Range("B37").Value = 5.20914
$wsDetailed.Range("B38").Value = 9.815239999999999
$wsDetailed.Range("B39").Value = 61.29848
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 67.22089
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 78
$wsDetailed.Range("B42").Value = 79.78207
$wsDetailed.Range("B44").Value = 77.94
$wsDetailed.Range("B45").Value = 65
$wsDetailed.Range("B47").Value = 64.8901
$wsDetailed.Range("B48").Value = 63.24225
$wsDetailed.Range("B59").Value = 68.64169
$wsDetailed.Range("B60").Value = 69.15218
$wsDetailed.Range("B61").Value = 79.95022
$wsDetailed.Range("B62").Value = 79.95016
$wsDetailed.Range("B63").Value = 63.13224
$wsDetailed.Range("B64").Value = 31.22136
$wsDetailed.Range("B65").Value = 0.6787
$wsDetailed.Range("B66").Value = -1.02778
$wsDetailed.Range("B67").Value = -5.74313
$wsDetailed.Range("B68").Value = -6.99878
$wsDetailed.Range("B69").Value = -7.12619
$wsDetailed.Range("B70").Value = -9.093769999999999
$wsDetailed.Range("B73").Value = -14.15388
$wsDetailed.Range("B74").Value = -14.51348
$wsDetailed.Range("B75").Value = -20
$wsDetailed.Range("B76").Value = -19.85975
$wsDetailed.Range("B77").Value = -23.03945
$wsDetailed.Range("B78").Value = -23.5
$wsDetailed.Range("B79").Value = -20.47588
$wsDetailed.Range("B80").Value = -22.88324
$wsDetailed.Range("B81").Value = -14
$wsDetailed.Range("B82").Value = -6.71274
$wsDetailed.Range("B83").Value = -5.01
$wsDetailed.Range("B85").Value = 42.24206
$wsDetailed.Range("B86").Value = 56.51941
$wsDetailed.Range("B87").Value = 65
$wsDetailed.Range("B89").Value = 103.6301
$wsDetailed.Range("B90").Value = 86.52197
$wsDetailed.Range("B91").Value = 73.65000000000001
$wsDetailed.Range("B92").Value = 67.87605000000001
$wsDetailed.Range("B94").Value = 59.90651
$wsDetailed.Range("B95").Value = 64.8901
$wsDetailed.Range("B96").Value = 64.23224
$wsDetailed.Range("B97").Value = 63.23644
